$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell 'D2' '57.964.74'
Set-TextCell 'E2' '  -3.15%  '
Set-TextCell 'D3' '2.284.93'
Set-TextCell 'E3' '  -3.76%  '
Set-TextCell 'D5' '533.52'
Set-TextCell 'E5' '  -4.14%  '
Set-TextCell 'E6' '  -1.75%  '
Set-TextCell 'E7' '  +0.05%  '
Set-TextCell 'D8' '0.588'
Set-TextCell 'E8' '  +0.73%  '
Set-TextCell 'D9' '2.283.87'
Set-TextCell 'E9' '  -3.75%  '
Set-TextCell 'E10' '  -5.49%  '
Set-TextCell 'D11' '5.44'
Set-TextCell 'E13' '  -3.55%  '
Set-TextCell 'D14' '23.52'
Set-TextCell 'E14' '  -3.55%  '
Set-TextCell 'D15' '2.694.19'
Set-TextCell 'E15' '  -3.73%  '
Set-TextCell 'D16' '57.913.57'
Set-TextCell 'E16' '  -3.16%  '
Set-TextCell 'E17' '  -4.41%  '
Set-TextCell 'D18' '2.293.67'
Set-TextCell 'E18' '  -3.22%  '
Set-TextCell 'D19' '10.53'
Set-TextCell 'E19' '  -5.19%  '
Set-TextCell 'E20' '  -5.60%  '
Set-TextCell 'D21' '312.86'
Set-TextCell 'E21' '  -2.33%  '
Set-TextCell 'D22' '6.37'
Set-TextCell 'E22' '  -4.31%  '
Set-TextCell 'E23' '  +0.09%  '
Set-TextCell 'D24' '62.42'
Set-TextCell 'E24' '  -2.76%  '
Set-TextCell 'E25' '  -3.30%  '
Set-TextCell 'E26' '  -0.07%  '
Set-TextCell 'E27' '  -5.38%  '
Set-TextCell 'E28' '  -7.28%  '
Set-TextCell 'D29' '170.87'
Set-TextCell 'E29' '  +0.36%  '
Set-TextCell 'E30' '  -5.73%  '
Set-TextCell 'E31' '  -5.26%  '
Set-TextCell 'D33' '1.05'
Set-TextCell 'E33' '  -6.85%  '
Set-TextCell 'D34' '0.379'
Set-TextCell 'E34' '  -5.68%  '
Set-TextCell 'E36' '  -2.16%  '
Set-TextCell 'E37' '  -0.08%  '
Set-TextCell 'E38' '  -7.48%  '
Set-TextCell 'E39' '  -5.93%  '
Set-TextCell 'E40' '  -1.21%  '
Set-TextCell 'D41' '1.49'
Set-TextCell 'E41' '  -6.43%  '
Set-TextCell 'D42' '141.82'
Set-TextCell 'E42' '  -2.31%  '
Set-TextCell 'D43' '286.73'
Set-TextCell 'E43' '  -10.57%  '
Set-TextCell 'E44' '  -2.97%  '
Set-TextCell 'D45' '0.0952'
Set-TextCell 'E45' '  -1.43%  '
Set-TextCell 'D46' '0.0495'
Set-TextCell 'E46' '  -3.10%  '
Set-TextCell 'D47' '0.552'
Set-TextCell 'E47' '  -2.80%  '
Set-TextCell 'D48' '18.05'
Set-TextCell 'E48' '  -8.08%  '
Set-TextCell 'D49' '0.0211'
Set-TextCell 'E49' '  -3.17%  '
Set-TextCell 'D50' '10.92'
Set-TextCell 'E50' '  -1.31%  '
Set-TextCell 'E51' '  -0.70%  '
